# Stevens & Palocsay (2017) mean-vs-variance example:
# add a "Supervisors" sheet, and link each Project to a Supervisor.

$wb = $excel.ActiveWorkbook

# --- New "Supervisors" sheet (inserted as the very first sheet) -----------
$supervisors = $wb.Worksheets.Add()
$supervisors.Name = "Supervisors"

$supervisors.Range("A1").Value = "Supervisor"
$supervisors.Range("B1").Value = "Max_number_of_projects"
$supervisors.Range("C1").Value = "Max_number_of_students"
$supervisors.Range("A2").Value = "Dr Smith"

$supervisors.Columns.Item(2).ColumnWidth = 21.8

# --- "Projects" sheet: rename header & add a Supervisor column ------------
$projects = $wb.Worksheets.Item("Projects")

$projects.Range("A1").Value = "Project"
$projects.Range("C1").Value = "Supervisor"

$projects.Range("C2").Value = "Dr Smith"
$projects.Range("C3").Value = "Dr Smith"
$projects.Range("C4").Value = "Dr Smith"
$projects.Range("C5").Value = "Dr Smith"
$projects.Range("C6").Value = "Dr Smith"

$projects.Columns.Item(2).ColumnWidth = 22.2

# --- Restore "Projects" as the active / selected sheet ---------------------
$projects.Activate()
